$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection to G2 (matches sheetView selection change E6 -> G2)
$ws.Range("G2").Select()

# Update the "Date of Execution" for the existing "Generate report" row (row 2)
$ws.Range("G2").Value = "1/17/2017"

# Fill in the new "Save Report" sprint task row (row 3)
$ws.Range("A3").Value = "Save Report"
$ws.Range("B3").Value = "T_C_"
$ws.Range("C3").Value = "Save Report"
$ws.Range("E3").Value = "Medium"
$ws.Range("G3").Value = "1/17/2017"
